# Auto-generated edit script: updates market-price-derived cells per the commit diff.
# Applies changed/added/removed numeric cell values across the 8 Leve-profit tables (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 2325.375
$ws.Range("I12").Value = 2520
$ws.Range("K12").Value = 2520
$ws.Range("M12").Value = -2350
# Row 40
$ws.Range("I40").Value = 2850
$ws.Range("J40").Value = 4038.2354
$ws.Range("K40").Value = 2850
$ws.Range("L40").Value = 4038.2354
$ws.Range("M40").Value = -2675
$ws.Range("N40").Value = -4388.2354
# Row 92
$ws.Range("H92").Value = 1205.8823
$ws.Range("I92").Value = 774.6
$ws.Range("K92").Value = 774.6
$ws.Range("M92").Value = 473.4
# Row 101
$ws.Range("H101").Value = 2503.4736
$ws.Range("I101").Value = 1423.7273
$ws.Range("J101").Value = 3988.125
$ws.Range("K101").Value = 4271.1819
$ws.Range("L101").Value = 11964.375
$ws.Range("M101").Value = -2649.1819
$ws.Range("N101").Value = -15208.375
# Row 116
$ws.Range("H116").Value = 6180.1875
$ws.Range("I116").Value = 6043.222
$ws.Range("K116").Value = 6043.222
$ws.Range("M116").Value = -2601.222
# Row 138
$ws.Range("H138").Value = 4680.7646
$ws.Range("I138").Value = 4469.5
$ws.Range("K138").Value = 13408.5
$ws.Range("M138").Value = -8268.5
# Row 141
$ws.Range("H141").Value = 1627
$ws.Range("I141").Value = 1474.8572
$ws.Range("J141").Value = 2159.5
$ws.Range("K141").Value = 4424.571599999999
$ws.Range("L141").Value = 6478.5
$ws.Range("M141").Value = 755.4284000000007
$ws.Range("N141").Value = -16838.5

$ws = $wb.Worksheets.Item("ARM")
# Row 13
$ws.Range("H13").Value = 21950.166
$ws.Range("I13").Value = 1001
$ws.Range("J13").Value = 26140
$ws.Range("K13").Value = 1001
$ws.Range("L13").Value = 26140
$ws.Range("N13").Value = -26428
$ws.Range("M13").Value = -857
# Row 43
$ws.Range("H43").Value = 27450
$ws.Range("J43").Value = 27450
$ws.Range("L43").Value = 27450
$ws.Range("N43").Value = -28076
# Row 45
$ws.Range("H45").Value = 2275.25
$ws.Range("I45").Value = 1834.6666
$ws.Range("J45").Value = 2464.0715
$ws.Range("K45").Value = 1834.6666
$ws.Range("L45").Value = 2464.0715
$ws.Range("M45").Value = -1457.6666
$ws.Range("N45").Value = -3218.0715
# Row 61
$ws.Range("H61").Value = 41761380
$ws.Range("I61").Value = 166675000
$ws.Range("J61").Value = 123507.445
$ws.Range("K61").Value = 166675000
$ws.Range("L61").Value = 123507.445
$ws.Range("M61").Value = -166674788
$ws.Range("N61").Value = -123931.445
# Row 74
$ws.Range("H74").Value = 6950340.5
$ws.Range("I74").Value = 12501230
$ws.Range("J74").Value = 11729.625
$ws.Range("K74").Value = 12501230
$ws.Range("L74").Value = 11729.625
$ws.Range("M74").Value = -12500356
$ws.Range("N74").Value = -13477.625
# Row 77
$ws.Range("H77").Value = 6950340.5
$ws.Range("I77").Value = 12501230
$ws.Range("J77").Value = 11729.625
$ws.Range("K77").Value = 62506150
$ws.Range("L77").Value = 58648.125
$ws.Range("M77").Value = -62501782
$ws.Range("N77").Value = -67384.125
# Row 97
$ws.Range("H97").Value = 1209.6786
$ws.Range("I97").Value = 791.5217
$ws.Range("K97").Value = 791.5217
$ws.Range("M97").Value = -295.5217
# Row 122
$ws.Range("H122").Value = 2878.3928
$ws.Range("I122").Value = 1701.4445
$ws.Range("K122").Value = 5104.333500000001
$ws.Range("M122").Value = -2654.333500000001
# Row 132
$ws.Range("H132").Value = 8007.775
$ws.Range("I132").Value = 4645.0356
$ws.Range("K132").Value = 13935.1068
$ws.Range("M132").Value = -11405.1068
# Row 136
$ws.Range("H136").Value = 41761380
$ws.Range("I136").Value = 166675000
$ws.Range("J136").Value = 123507.445
$ws.Range("K136").Value = 500025000
$ws.Range("L136").Value = 370522.335
$ws.Range("M136").Value = -500022450
$ws.Range("N136").Value = -375622.335
# Row 138
$ws.Range("H138").Value = 275000
$ws.Range("J138").Value = 275000
$ws.Range("L138").Value = 275000
$ws.Range("N138").Value = -285280

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 4857.5293
$ws.Range("I86").Value = 3564.2144
$ws.Range("K86").Value = 3564.2144
$ws.Range("M86").Value = -2441.2144
# Row 89
$ws.Range("H89").Value = 4857.5293
$ws.Range("I89").Value = 3564.2144
$ws.Range("K89").Value = 17821.072
$ws.Range("M89").Value = -12205.072
# Row 134
$ws.Range("H134").Value = 25741.861
$ws.Range("I134").Value = 2719.5667
$ws.Range("J134").Value = 78870.234
$ws.Range("K134").Value = 8158.7001
$ws.Range("L134").Value = 236610.702
$ws.Range("M134").Value = -5623.7001
$ws.Range("N134").Value = -241680.702

$ws = $wb.Worksheets.Item("CRP")
# Row 3
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
# Row 7
$ws.Range("H7").Value = 191.125
$ws.Range("I7").Value = 171.6
$ws.Range("J7").Value = 223.66667
$ws.Range("K7").Value = 171.6
$ws.Range("L7").Value = 223.66667
$ws.Range("M7").Value = -58.59999999999999
$ws.Range("N7").Value = -449.66667
# Row 31
$ws.Range("H31").Value = 1306519.5
$ws.Range("I31").Value = 27917
$ws.Range("K31").Value = 27917
$ws.Range("M31").Value = -27622
# Row 34
$ws.Range("H34").Value = 1306519.5
$ws.Range("I34").Value = 27917
$ws.Range("K34").Value = 27917
$ws.Range("M34").Value = -27715
# Row 57
$ws.Range("H57").Value = 20000
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 20000
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 20000
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -21120
# Row 63
$ws.Range("H63").Value = 48835.5
$ws.Range("J63").Value = 48835.5
$ws.Range("L63").Value = 48835.5
$ws.Range("N63").Value = -50207.5
# Row 66
$ws.Range("H66").Value = 48835.5
$ws.Range("J66").Value = 48835.5
$ws.Range("L66").Value = 146506.5
$ws.Range("N66").Value = -153370.5
# Row 134
$ws.Range("H134").Value = 719815.6
$ws.Range("J134").Value = 15246.75
$ws.Range("L134").Value = 45740.25
$ws.Range("N134").Value = -50810.25

$ws = $wb.Worksheets.Item("CUL")
# Row 8
$ws.Range("H8").Value = 311.33334
$ws.Range("I8").Value = 311.33334
$ws.Range("K8").Value = 934.0000200000001
$ws.Range("M8").Value = -795.0000200000001
# Row 60
$ws.Range("H60").Value = 1346.3334
$ws.Range("I60").Value = 1052.5
$ws.Range("K60").Value = 3157.5
$ws.Range("M60").Value = -2906.5
# Row 131
$ws.Range("H131").Value = 3706.65
$ws.Range("I131").Value = 3914
$ws.Range("J131").Value = 2531.6667
$ws.Range("K131").Value = 11742
$ws.Range("L131").Value = 7595.000100000001
$ws.Range("M131").Value = -6702
$ws.Range("N131").Value = -17675.0001

$ws = $wb.Worksheets.Item("GSM")
# Row 31
$ws.Range("H31").Value = 3563.3333
$ws.Range("I31").Value = 3563.3333
$ws.Range("K31").Value = 3563.3333
$ws.Range("M31").Value = -3271.3333
# Row 37
$ws.Range("H37").Value = 3563.3333
$ws.Range("I37").Value = 3563.3333
$ws.Range("K37").Value = 3563.3333
$ws.Range("M37").Value = -3286.3333
# Row 80
$ws.Range("H80").Value = 12767.177
$ws.Range("I80").Value = 5671
$ws.Range("K80").Value = 5671
$ws.Range("M80").Value = -4673
# Row 83
$ws.Range("H83").Value = 12767.177
$ws.Range("I83").Value = 5671
$ws.Range("K83").Value = 28355
$ws.Range("M83").Value = -23363

$ws = $wb.Worksheets.Item("LTW")
# Row 12
$ws.Range("H12").Value = 20500.75
$ws.Range("I12").Value = 1333
$ws.Range("J12").Value = 32001.4
$ws.Range("K12").Value = 1333
$ws.Range("L12").Value = 32001.4
$ws.Range("M12").Value = -1163
$ws.Range("N12").Value = -32341.4
# Row 22
$ws.Range("H22").Value = 2994.7273
$ws.Range("I22").Value = 2797.9666
$ws.Range("K22").Value = 2797.9666
$ws.Range("M22").Value = -2502.9666
# Row 27
$ws.Range("H27").Value = 2994.7273
$ws.Range("I27").Value = 2797.9666
$ws.Range("K27").Value = 2797.9666
$ws.Range("M27").Value = -2690.9666
# Row 46
$ws.Range("H46").Value = 1987.9445
$ws.Range("I46").Value = 1832.5555
$ws.Range("K46").Value = 1832.5555
$ws.Range("M46").Value = -1644.5555
# Row 55
$ws.Range("H55").Value = 100000720
$ws.Range("I55").Value = 125000670
$ws.Range("K55").Value = 125000670
$ws.Range("M55").Value = -125000497
# Row 70
$ws.Range("H70").Value = 54163
$ws.Range("J70").Value = 54163
$ws.Range("L70").Value = 54163
$ws.Range("N70").Value = -54703
# Row 73
$ws.Range("H73").Value = 54163
$ws.Range("J73").Value = 54163
$ws.Range("L73").Value = 54163
$ws.Range("N73").Value = -56035

$ws = $wb.Worksheets.Item("WVR")
# Row 7
$ws.Range("H7").Value = 19286.428
$ws.Range("I7").Value = 20000
$ws.Range("K7").Value = 20000
$ws.Range("M7").Value = -19887
# Row 14
$ws.Range("H14").Value = 2694.2
$ws.Range("I14").Value = 2694.2
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 2694.2
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -2526.2
$ws.Range("N14").ClearContents()
# Row 63
$ws.Range("H63").Value = 19332.334
$ws.Range("J63").Value = 17998.5
$ws.Range("L63").Value = 17998.5
$ws.Range("N63").Value = -19246.5
# Row 64
$ws.Range("H64").Value = 64900
$ws.Range("J64").Value = 64900
$ws.Range("L64").Value = 64900
$ws.Range("N64").Value = -65396
# Row 66
$ws.Range("H66").Value = 19332.334
$ws.Range("J66").Value = 17998.5
$ws.Range("L66").Value = 53995.5
$ws.Range("N66").Value = -60235.5
# Row 67
$ws.Range("H67").Value = 64900
$ws.Range("J67").Value = 64900
$ws.Range("L67").Value = 64900
$ws.Range("N67").Value = -66616
# Row 136
$ws.Range("H136").Value = 2952.6155
$ws.Range("I136").Value = 2018.5883
$ws.Range("K136").Value = 6055.7649
$ws.Range("M136").Value = -3505.7649
# Row 137
$ws.Range("H137").Value = 106355
$ws.Range("J137").Value = 106355
$ws.Range("L137").Value = 106355
$ws.Range("N137").Value = -116555
